$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personal")

# Stundensatz (hourly rate) gesenkt von 120 auf 80
$ws.Range("E1").Value = 80

# Selektierte Zelle passend zum Ziel-Dokument setzen
$ws.Range("F7").Select()

$wb.Save()
